# Fruta / hortaliza, semanal
# Insert two new weekly price rows (Carson - Primera / Segunda) for
# Vega Monumental Concepción - Durazno, pushing the existing data down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 226, shifting all
# subsequent rows (previously 226-283) down to 228-285.
$ws.Rows.Item(226).Resize(2).Insert()

# New row 226: Carson - Primera
$ws.Range("A226").Value = 11
$ws.Range("B226").Value = "Vega Monumental Concepción"
$ws.Range("C226").Value = "Bíobío"
$ws.Range("D226").Value = 44951
$ws.Range("E226").Value = 8
$ws.Range("F226").Value = "Fruta"
$ws.Range("G226").Value = 100103
$ws.Range("H226").Value = "Frutos de hueso (carozo)"
$ws.Range("I226").Value = 100103004
$ws.Range("J226").Value = "Durazno"
$ws.Range("K226").Value = "Carson"
$ws.Range("L226").Value = "Primera"
$ws.Range("M226").Value = 150
$ws.Range("N226").Value = 13000
$ws.Range("O226").Value = 13000
$ws.Range("P226").Value = 13000
$ws.Range("Q226").Value = "`$/caja 16 kilos empedrada"
$ws.Range("R226").Value = "Región de O'Higgins"
$ws.Range("S226").Value = 812
$ws.Range("T226").Value = 16

# New row 227: Carson - Segunda
$ws.Range("A227").Value = 11
$ws.Range("B227").Value = "Vega Monumental Concepción"
$ws.Range("C227").Value = "Bíobío"
$ws.Range("D227").Value = 44951
$ws.Range("E227").Value = 8
$ws.Range("F227").Value = "Fruta"
$ws.Range("G227").Value = 100103
$ws.Range("H227").Value = "Frutos de hueso (carozo)"
$ws.Range("I227").Value = 100103004
$ws.Range("J227").Value = "Durazno"
$ws.Range("K227").Value = "Carson"
$ws.Range("L227").Value = "Segunda"
$ws.Range("M227").Value = 100
$ws.Range("N227").Value = 10000
$ws.Range("O227").Value = 10000
$ws.Range("P227").Value = 10000
$ws.Range("Q227").Value = "`$/caja 16 kilos empedrada"
$ws.Range("R227").Value = "Región de O'Higgins"
$ws.Range("S227").Value = 625
$ws.Range("T227").Value = 16
